$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label for the IP address column
$ws.Range("B1").Value = "IP Address"

# Reflect the new selection/active cell recorded in the sheet view
$ws.Range("D6").Select()
